$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 17 and row 18 (columns B:AD), keep column A fixed
$ws.Range("B17").Value = 6838980
$ws.Range("B18").Value = 6838982
$ws.Range("E17").Value = 'Al Karkh'
$ws.Range("E18").Value = 'Al Sinaah'
$ws.Range("F17").Value = 'Al Talaba'
$ws.Range("F18").Value = 'Erbil SC'
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 2
$ws.Range("J17").Value = 0
$ws.Range("J18").Value = 2
$ws.Range("K17").Value = 'H'
$ws.Range("K18").Value = 'A'
$ws.Range("L17").Value = 3.2
$ws.Range("L18").Value = 3.6
$ws.Range("M17").Value = 2.8
$ws.Range("M18").Value = 3
$ws.Range("N17").Value = 2.25
$ws.Range("N18").Value = 2
$ws.Range("O17").Value = 3.1
$ws.Range("O18").Value = 3.5
$ws.Range("P17").Value = 2.8
$ws.Range("P18").Value = 3
$ws.Range("Q17").Value = 2.25
$ws.Range("Q18").Value = 2
$ws.Range("R17").Value = 0.25
$ws.Range("R18").Value = 0.5
$ws.Range("S17").Value = 1.8
$ws.Range("S18").Value = 1.725
$ws.Range("T17").Value = 2
$ws.Range("T18").Value = 1.975
$ws.Range("V17").Value = 1.85
$ws.Range("V18").Value = 2
$ws.Range("W17").Value = 1.95
$ws.Range("W18").Value = 1.8
$ws.Range("X17").Value = 2.1
$ws.Range("X18").Value = -1
$ws.Range("Z17").Value = -1
$ws.Range("Z18").Value = 1
$ws.Range("AA17").Value = 0.8
$ws.Range("AA18").Value = -1
$ws.Range("AB17").Value = -1
$ws.Range("AB18").Value = 0.9750000000000001
$ws.Range("AC17").Value = -1
$ws.Range("AC18").Value = 1
$ws.Range("AD17").Value = 0.95
$ws.Range("AD18").Value = -1

# Swap row 58 and row 59 (columns B:AD), keep column A fixed
$ws.Range("B58").Value = 7511489
$ws.Range("B59").Value = 7511490
$ws.Range("E58").Value = 'Naft AlBasra'
$ws.Range("E59").Value = 'Amanat Baghdad'
$ws.Range("F58").Value = 'Naft Maysan'
$ws.Range("F59").Value = 'Newroz SC'
$ws.Range("L58").Value = 2.2
$ws.Range("L59").Value = 3.75
$ws.Range("M58").Value = 3
$ws.Range("M59").Value = 3.1
$ws.Range("N58").Value = 3.1
$ws.Range("N59").Value = 1.909
$ws.Range("O58").Value = 2.2
$ws.Range("O59").Value = 2.875
$ws.Range("P58").Value = 3
$ws.Range("P59").Value = 2.9
$ws.Range("Q58").Value = 3.1
$ws.Range("Q59").Value = 2.3
$ws.Range("R58").Value = -0.25
$ws.Range("R59").Value = 0.25
$ws.Range("S58").Value = 1.95
$ws.Range("S59").Value = 1.75
$ws.Range("T58").Value = 1.85
$ws.Range("T59").Value = 2.05
$ws.Range("U58").Value = 2.25
$ws.Range("U59").Value = 1.75
$ws.Range("V58").Value = 2.05
$ws.Range("V59").Value = 1.9
$ws.Range("W58").Value = 1.75
$ws.Range("W59").Value = 1.9
$ws.Range("Z58").Value = 2.1
$ws.Range("Z59").Value = 1.3
$ws.Range("AB58").Value = 0.8500000000000001
$ws.Range("AB59").Value = 1.05
$ws.Range("AC58").Value = -0.5
$ws.Range("AC59").Value = 0.45
$ws.Range("AD58").Value = 0.375
$ws.Range("AD59").Value = -0.5

# Swap row 135 and row 136 (columns B:AD), keep column A fixed
$ws.Range("B135").Value = 7897090
$ws.Range("B136").Value = 7901411
$ws.Range("E135").Value = 'Al Karkh'
$ws.Range("E136").Value = 'Naft Maysan'
$ws.Range("F135").Value = 'Al Najaf'
$ws.Range("F136").Value = 'Al Naft SC'
$ws.Range("H135").Value = 3
$ws.Range("H136").Value = 0
$ws.Range("J135").Value = 1
$ws.Range("J136").Value = 0
$ws.Range("K135").Value = 'A'
$ws.Range("K136").Value = 'D'
$ws.Range("M135").Value = 2.8
$ws.Range("M136").Value = 2.9
$ws.Range("N135").Value = 3.6
$ws.Range("N136").Value = 3.4
$ws.Range("O135").Value = 2.1
$ws.Range("O136").Value = 2.15
$ws.Range("P135").Value = 2.8
$ws.Range("P136").Value = 2.875
$ws.Range("Q135").Value = 3.6
$ws.Range("Q136").Value = 3.3
$ws.Range("S135").Value = 1.9
$ws.Range("S136").Value = 1.925
$ws.Range("T135").Value = 1.9
$ws.Range("T136").Value = 1.875
$ws.Range("U135").Value = 2
$ws.Range("U136").Value = 1.75
$ws.Range("V135").Value = 2
$ws.Range("V136").Value = 1.75
$ws.Range("W135").Value = 1.8
$ws.Range("W136").Value = 2.05
$ws.Range("Y135").Value = -1
$ws.Range("Y136").Value = 1.875
$ws.Range("Z135").Value = 2.6
$ws.Range("Z136").Value = -1
$ws.Range("AA135").Value = -1
$ws.Range("AA136").Value = -0.5
$ws.Range("AB135").Value = 0.8999999999999999
$ws.Range("AB136").Value = 0.4375
$ws.Range("AC135").Value = 1
$ws.Range("AC136").Value = -1
$ws.Range("AD135").Value = -1
$ws.Range("AD136").Value = 1.05

# Swap row 173 and row 174 (columns B:AD), keep column A fixed
$ws.Range("B173").Value = 8137759
$ws.Range("B174").Value = 8137758
$ws.Range("E173").Value = 'Al Quwa Al Jawiya'
$ws.Range("E174").Value = 'Al Minaa'
$ws.Range("F173").Value = 'Karbalaa FC'
$ws.Range("F174").Value = 'Newroz SC'
$ws.Range("G173").Value = 2
$ws.Range("G174").Value = 1
$ws.Range("H173").Value = 1
$ws.Range("H174").Value = 4
$ws.Range("I173").Value = 1
$ws.Range("I174").Value = 0
$ws.Range("J173").Value = 0
$ws.Range("J174").Value = 2
$ws.Range("K173").Value = 'H'
$ws.Range("K174").Value = 'A'
$ws.Range("L173").Value = 1.571
$ws.Range("L174").Value = 4.333
$ws.Range("M173").Value = 3.6
$ws.Range("M174").Value = 3.8
$ws.Range("N173").Value = 5
$ws.Range("N174").Value = 1.615
$ws.Range("O173").Value = 1.25
$ws.Range("O174").Value = 3.6
$ws.Range("P173").Value = 4.5
$ws.Range("P174").Value = 3.5
$ws.Range("Q173").Value = 10
$ws.Range("Q174").Value = 1.8
$ws.Range("R173").Value = -1.75
$ws.Range("R174").Value = 0.5
$ws.Range("S173").Value = 2
$ws.Range("S174").Value = 1.975
$ws.Range("T173").Value = 1.8
$ws.Range("T174").Value = 1.825
$ws.Range("U173").Value = 2.75
$ws.Range("U174").Value = 2.25
$ws.Range("V173").Value = 1.925
$ws.Range("V174").Value = 1.975
$ws.Range("W173").Value = 1.875
$ws.Range("W174").Value = 1.725
$ws.Range("X173").Value = 0.25
$ws.Range("X174").Value = -1
$ws.Range("Z173").Value = -1
$ws.Range("Z174").Value = 0.8
$ws.Range("AB173").Value = 0.8
$ws.Range("AB174").Value = 0.825
$ws.Range("AC173").Value = 0.4625
$ws.Range("AC174").Value = 0.9750000000000001
$ws.Range("AD173").Value = -0.5
$ws.Range("AD174").Value = -1

# Swap row 221 and row 222 (columns B:AD), keep column A fixed
$ws.Range("B221").Value = 8259836
$ws.Range("B222").Value = 8261818
$ws.Range("E221").Value = 'Al Hudod'
$ws.Range("E222").Value = 'Newroz SC'
$ws.Range("F221").Value = 'Al Talaba'
$ws.Range("F222").Value = 'Naft Maysan'
$ws.Range("G221").Value = 1
$ws.Range("G222").Value = 0
$ws.Range("H221").Value = 2
$ws.Range("H222").Value = 0
$ws.Range("K221").Value = 'A'
$ws.Range("K222").Value = 'D'
$ws.Range("L221").Value = 2.25
$ws.Range("L222").Value = 2.2
$ws.Range("M221").Value = 2.875
$ws.Range("M222").Value = 2.75
$ws.Range("N221").Value = 3.1
$ws.Range("N222").Value = 3.4
$ws.Range("O221").Value = 2.45
$ws.Range("O222").Value = 2.2
$ws.Range("P221").Value = 2.875
$ws.Range("P222").Value = 2.8
$ws.Range("Q221").Value = 2.8
$ws.Range("Q222").Value = 3.4
$ws.Range("R221").Value = 0
$ws.Range("R222").Value = -0.25
$ws.Range("S221").Value = 1.775
$ws.Range("S222").Value = 1.9
$ws.Range("T221").Value = 2.025
$ws.Range("T222").Value = 1.9
$ws.Range("U221").Value = 2
$ws.Range("U222").Value = 2.25
$ws.Range("V221").Value = 1.975
$ws.Range("V222").Value = 1.8
$ws.Range("W221").Value = 1.825
$ws.Range("W222").Value = 2
$ws.Range("Y221").Value = -1
$ws.Range("Y222").Value = 1.8
$ws.Range("Z221").Value = 1.8
$ws.Range("Z222").Value = -1
$ws.Range("AA221").Value = -1
$ws.Range("AA222").Value = -0.5
$ws.Range("AB221").Value = 1.025
$ws.Range("AB222").Value = 0.45
$ws.Range("AC221").Value = 0.9750000000000001
$ws.Range("AC222").Value = -1
$ws.Range("AD221").Value = -1
$ws.Range("AD222").Value = 1

# Swap row 223 and row 224 (columns B:AD), keep column A fixed
$ws.Range("B223").Value = 8263267
$ws.Range("B224").Value = 8267113
$ws.Range("E223").Value = 'Zakho'
$ws.Range("E224").Value = 'Karbalaa FC'
$ws.Range("F223").Value = 'Al Naft SC'
$ws.Range("F224").Value = 'Al Najaf'
$ws.Range("G223").Value = 3
$ws.Range("G224").Value = 0
$ws.Range("K223").Value = 'H'
$ws.Range("K224").Value = 'D'
$ws.Range("L223").Value = 2.2
$ws.Range("L224").Value = 3.75
$ws.Range("M223").Value = 2.8
$ws.Range("M224").Value = 2.875
$ws.Range("N223").Value = 3.25
$ws.Range("N224").Value = 2
$ws.Range("O223").Value = 2.375
$ws.Range("O224").Value = 3.2
$ws.Range("P223").Value = 2.625
$ws.Range("P224").Value = 2.8
$ws.Range("Q223").Value = 3.1
$ws.Range("Q224").Value = 2.25
$ws.Range("R223").Value = -0.25
$ws.Range("R224").Value = 0.25
$ws.Range("S223").Value = 2.025
$ws.Range("S224").Value = 1.825
$ws.Range("T223").Value = 1.775
$ws.Range("T224").Value = 1.975
$ws.Range("V223").Value = 2.025
$ws.Range("V224").Value = 1.975
$ws.Range("W223").Value = 1.775
$ws.Range("W224").Value = 1.825
$ws.Range("X223").Value = 1.375
$ws.Range("X224").Value = -1
$ws.Range("Y223").Value = -1
$ws.Range("Y224").Value = 1.8
$ws.Range("AA223").Value = 1.025
$ws.Range("AA224").Value = 0.4125
$ws.Range("AB223").Value = -1
$ws.Range("AB224").Value = -0.5
$ws.Range("AC223").Value = 1.025
$ws.Range("AC224").Value = -1
$ws.Range("AD223").Value = -1
$ws.Range("AD224").Value = 0.825
